# Bulk-app translation workbook: rename icon/audio filepath columns to
# their localized (per-language) equivalents, and restore the "natural"
# active sheet / selection state that Excel saves when a user is looking
# at the Modules_and_forms sheet with F2 selected.

$wb = $excel.ActiveWorkbook

$wsModules = $wb.Worksheets.Item("Modules_and_forms")

# --- Rename the two header cells that now carry a language suffix ------
# E1: icon_filepath  -> icon_filepath_en
# F1: audio_filepath -> audio_filepath_en
$wsModules.Cells.Item(1, 5).Value = "icon_filepath_en"
$wsModules.Cells.Item(1, 6).Value = "audio_filepath_en"

# --- Restore saved view state ------------------------------------------
# The workbook was last saved with "Modules_and_forms" active and cell
# F2 selected (instead of module1 / B11 being active).
$wsModules.Activate()
$wsModules.Range("F2").Select()
